# Generate Report for Handoff
#
# Adds two new handed-off files (99f2b05f-... and ac35203f-...) to the
# localization status report: one new row on "Overview", and one new row
# on each of the "zh-cn" and "de-de" language sheets. Resizes the backing
# Excel tables / sheet dimensions to match.

$wb = $excel.ActiveWorkbook

$dateDeDe = "2016-08-25 06:39:31"
$dateZhCn = "2016-08-25 06:39:26"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 4; Name = "99f2b05f-ea50-4f8e-8b54-98be03107cb1.md" },
    @{ Row = 5; Name = "ac35203f-9a7b-4797-baa8-8a9b214211cd.md" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $name = $r.Name
    $display = "e2e\" + $name

    $wsOverview.Range("A" + $row).Value = $name
    $wsOverview.Hyperlinks.Add($wsOverview.Range("B" + $row), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/" + $name, "", "", $display)
    $wsOverview.Range("C" + $row).Value = ".md"
    $wsOverview.Range("D" + $row).Value = "'"
    $wsOverview.Range("E" + $row).Value = "Ready for handoff"
    $wsOverview.Range("F" + $row).Value = "Ready for handoff"
    $wsOverview.Range("G" + $row).Value = $dateDeDe
    $wsOverview.Range("G" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------
# Language sheets (zh-cn / de-de)
# ---------------------------------------------------------------------
$languages = @(
    @{
        SheetName = "zh-cn";
        Date = $dateZhCn;
        Files = @(
            @{ Row = 4; Name = "99f2b05f-ea50-4f8e-8b54-98be03107cb1.md"; Xlf = "99f2b05f-ea50-4f8e-8b54-98be03107cb1.f1cb033a560a8c7cedc67c9a260d2326309238f6.zh-cn.xlf" },
            @{ Row = 5; Name = "ac35203f-9a7b-4797-baa8-8a9b214211cd.md"; Xlf = "ac35203f-9a7b-4797-baa8-8a9b214211cd.f1c75a17c9a06b18edc05e530f2330c797f8afd6.zh-cn.xlf" }
        )
    },
    @{
        SheetName = "de-de";
        Date = $dateDeDe;
        Files = @(
            @{ Row = 4; Name = "99f2b05f-ea50-4f8e-8b54-98be03107cb1.md"; Xlf = "99f2b05f-ea50-4f8e-8b54-98be03107cb1.f1cb033a560a8c7cedc67c9a260d2326309238f6.de-de.xlf" },
            @{ Row = 5; Name = "ac35203f-9a7b-4797-baa8-8a9b214211cd.md"; Xlf = "ac35203f-9a7b-4797-baa8-8a9b214211cd.f1c75a17c9a06b18edc05e530f2330c797f8afd6.de-de.xlf" }
        )
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.SheetName)

    foreach ($f in $lang.Files) {
        $row = $f.Row
        $name = $f.Name

        $ws.Hyperlinks.Add($ws.Range("A" + $row), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/" + $name, "", "", $name)
        $ws.Range("B" + $row).Value = ".md"
        $ws.Range("C" + $row).Value = "Ready for handoff"
        $ws.Range("D" + $row).Value = "e2e"
        $ws.Range("E" + $row).Value = "ht"
        $ws.Range("F" + $row).Value = "False"
        $ws.Range("G" + $row).Value = $f.Xlf
        $ws.Range("H" + $row).Value = $lang.Date
        $ws.Range("H" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("I" + $row).Value = ""
        $ws.Range("J" + $row).Value = ""
        $ws.Range("K" + $row).Value = "0001-01-01 00:00:00"
        $ws.Range("K" + $row).NumberFormat = "yyyy-mm-dd HH:mm:ss"
        $ws.Range("L" + $row).Value = ""
        $ws.Range("M" + $row).Value = "True"
        $ws.Range("N" + $row).Value = ""
        $ws.Range("O" + $row).Value = "False"
        $ws.Range("P" + $row).Value = ""
    }

    $ws.ListObjects.Item(1).Resize($ws.Range("A1:P5"))
}
